$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 187
$ws.Range("F6").Value = 5215
$ws.Range("F7").Value = 436
$ws.Range("F8").Value = 614
$ws.Range("F9").Value = 903
$ws.Range("F10").Value = 812
$ws.Range("F13").Value = 560
$ws.Range("F17").Value = 1739
$ws.Range("F18").Value = 1446
$ws.Range("F19").Value = 814
$ws.Range("F20").Value = 296
$ws.Range("F21").Value = 184
$ws.Range("F22").Value = 303
$ws.Range("F23").Value = 503
$ws.Range("F24").Value = 129
$ws.Range("F25").Value = 1044
$ws.Range("F27").Value = 525
$ws.Range("F28").Value = 2488
$ws.Range("F30").Value = 94
$ws.Range("F31").Value = 47
$ws.Range("F32").Value = 86
$ws.Range("F33").Value = 21
$ws.Range("F34").Value = 246
$ws.Range("F35").Value = 10
$ws.Range("F39").Value = 271
$ws.Range("F40").Value = 632
$ws.Range("F43").Value = 43
$ws.Range("F44").Value = 57

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 149
$ws.Range("F11").Value = 4

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 229

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 229
$ws.Range("F3").Value = 187
$ws.Range("F7").Value = 5215
$ws.Range("F8").Value = 436
$ws.Range("F9").Value = 614
$ws.Range("F11").Value = 149
$ws.Range("F12").Value = 903
$ws.Range("F13").Value = 812
$ws.Range("F17").Value = 560
$ws.Range("F22").Value = 1739
$ws.Range("F23").Value = 1446
$ws.Range("F24").Value = 814
$ws.Range("F25").Value = 296
$ws.Range("F26").Value = 184
$ws.Range("F27").Value = 303
$ws.Range("F29").Value = 503
$ws.Range("F30").Value = 129
$ws.Range("F31").Value = 1044
$ws.Range("F32").Value = 525
$ws.Range("F33").Value = 2489
$ws.Range("F35").Value = 94
$ws.Range("F36").Value = 86
$ws.Range("F37").Value = 21
$ws.Range("F38").Value = 246
$ws.Range("F39").Value = 10
$ws.Range("F42").Value = 271
$ws.Range("F43").Value = 632
$ws.Range("F45").Value = 43
$ws.Range("F46").Value = 57
$ws.Range("F49").Value = 4
